# Generate Report for Handback
# Updates the CI handback-status report with freshly regenerated timestamps
# (and a priority change) coming from a re-run of the handback report
# generator. The underlying shared-string values are updated everywhere
# they are used across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# "Latest HO Xliff Generate Date" for the 8922c5a8... entry (and the
# c9cea973... entry, which previously happened to share the same value)
$wsOverview.Range("G3").Value = "2016-08-17 12:15:14"
$wsOverview.Range("G4").Value = "2016-08-17 12:15:14"

# --- zh-cn sheet ------------------------------------------------------
# Priority changed from "ht" to "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# Correspond Handoff Datetime
$wsZhCn.Range("H3").Value = "2016-08-17 12:15:04"
$wsZhCn.Range("H4").Value = "2016-08-17 12:15:04"

# Correspond Handback DateTime
$wsZhCn.Range("K3").Value = "2016-08-17 12:15:29"
$wsZhCn.Range("K4").Value = "2016-08-17 12:15:29"

# --- de-de sheet --------------------------------------------------
# Priority changed from "ht" to "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# Correspond Handoff Datetime (re-uses the same value as Overview!G3/G4)
$wsDeDe.Range("H3").Value = "2016-08-17 12:15:14"
$wsDeDe.Range("H4").Value = "2016-08-17 12:15:14"

# Correspond Handback DateTime
$wsDeDe.Range("K3").Value = "2016-08-17 12:15:37"
$wsDeDe.Range("K4").Value = "2016-08-17 12:15:37"
